# New crime data collected - weekly CompStat 104th Precinct update
# (Volume/Number header + reporting week dates, then the Week-to-Date /
#  28-Day / Year-to-Date / 2-Year crime complaint stats table.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header text: "Volume 32   Number  27" -> "...Number  28" ----
$ws.Range("A8").Value = "Volume 32   Number  28"

# ---- Header text: reporting week dates ----
$ws.Range("C9").Value = "Report Covering the Week  7/7/2025  Through  7/13/2025"

# ---- Row 15 (Rape) ----
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 100
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = -5.882352941176

# ---- Row 16 (Robbery) ----
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 3
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = -11.764705882352
$ws.Range("I16").Value = 74
$ws.Range("J16").Value = 113
$ws.Range("K16").Value = -34.513274336283
$ws.Range("L16").Value = -36.206896551724
$ws.Range("M16").Value = -43.076923076923
$ws.Range("N16").Value = -85.346534653465

# ---- Row 17 (Fel. Assault) ----
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 42.857142857142
$ws.Range("F17").Value = 29
$ws.Range("G17").Value = 29
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 172
$ws.Range("J17").Value = 171
$ws.Range("K17").Value = 0.584795321637
$ws.Range("L17").Value = 10.256410256410
$ws.Range("M17").Value = 38.709677419354
$ws.Range("N17").Value = 12.418300653594

# ---- Row 18 (Burglary) ----
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = -37.5
$ws.Range("I18").Value = 88
$ws.Range("J18").Value = 104
$ws.Range("K18").Value = -15.384615384615
$ws.Range("L18").Value = -8.333333333333
$ws.Range("M18").Value = -63.333333333333
$ws.Range("N18").Value = -91.934005499541

# ---- Row 19 (Gr. Larceny) ----
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -13.333333333333
$ws.Range("F19").Value = 39
$ws.Range("G19").Value = 53
$ws.Range("H19").Value = -26.415094339622
$ws.Range("I19").Value = 321
$ws.Range("J19").Value = 372
$ws.Range("K19").Value = -13.709677419354
$ws.Range("L19").Value = -8.022922636103
$ws.Range("M19").Value = 47.247706422018
$ws.Range("N19").Value = -1.230769230769

# ---- Row 20 (G.L.A.) ----
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 12
$ws.Range("E20").Value = -25
$ws.Range("F20").Value = 32
$ws.Range("G20").Value = 43
$ws.Range("H20").Value = -25.581395348837
$ws.Range("I20").Value = 166
$ws.Range("J20").Value = 205
$ws.Range("K20").Value = -19.024390243902
$ws.Range("L20").Value = -6.214689265536
$ws.Range("M20").Value = -15.736040609137
$ws.Range("N20").Value = -91.327063740856

# ---- Row 21 (TOTAL) ----
$ws.Range("C21").Value = 41
$ws.Range("D21").Value = 41
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 128
$ws.Range("G21").Value = 159
$ws.Range("H21").Value = -19.496855345911
$ws.Range("I21").Value = 839
$ws.Range("J21").Value = 982
$ws.Range("K21").Value = -14.562118126272
$ws.Range("L21").Value = -7.497243660418
$ws.Range("M21").Value = -8.605664488017
$ws.Range("N21").Value = -79.103362391033

# ---- Row 22 (Transit) ----
# C22 switches from the "n/a" text placeholder (shared string "0") to an
# actual numeric count of 1 - copy number formatting from a same-style
# numeric neighbour (F22) after setting the value so the stored cell type
# flips from string to number.
$ws.Range("C22").Value = 1
$ws.Range("F22").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 10
$ws.Range("K22").Value = -23.076923076923
$ws.Range("L22").Value = -28.571428571428
$ws.Range("M22").Value = 0

# ---- Row 24 (Petit Larceny) ----
$ws.Range("C24").Value = 30
$ws.Range("E24").Value = -6.25
$ws.Range("F24").Value = 95
$ws.Range("G24").Value = 124
$ws.Range("H24").Value = -23.387096774193
$ws.Range("I24").Value = 644
$ws.Range("J24").Value = 768
$ws.Range("K24").Value = -16.145833333333
$ws.Range("L24").Value = -4.592592592592
$ws.Range("M24").Value = 11.226252158894

# ---- Row 25 (Retail Theft) ----
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = -42.857142857142
$ws.Range("G25").Value = 50
$ws.Range("H25").Value = -44
$ws.Range("I25").Value = 252
$ws.Range("J25").Value = 279
$ws.Range("K25").Value = -9.677419354838
$ws.Range("L25").Value = 12

# ---- Row 26 (Misd. Assault) ----
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -20
$ws.Range("F26").Value = 45
$ws.Range("G26").Value = 53
$ws.Range("H26").Value = -15.094339622641
$ws.Range("I26").Value = 290
$ws.Range("J26").Value = 333
$ws.Range("K26").Value = -12.912912912912
$ws.Range("L26").Value = 8.208955223880
$ws.Range("M26").Value = -28.395061728395

# ---- Row 27 (UCR Rape*) ----
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 100
$ws.Range("L27").Value = 0

# ---- Row 28 (Other Sex Crimes) ----
# F28 switches from a numeric count of 2 to the "n/a" text placeholder
# (shared string "0") - force text storage via the Text number format,
# then restore the normal formatting for this row by copying format
# from a same-style text neighbour (C14).
$ws.Range("F28").NumberFormat = "@"
$ws.Range("F28").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F28").PasteSpecial(-4122)

$ws.Range("G28").Value = 1
$ws.Range("H28").Value = -100
$ws.Range("L28").Value = -18.421052631578

# ---- Row 29 (Shooting Vic.) ----
$ws.Range("N29").Value = -71.428571428571

# ---- Row 30 (Shooting Inc.) ----
$ws.Range("N30").Value = -76.923076923076
